# Weekly update: insert a new data row at row 76 (pushing existing rows 76-188 down to 77-189)
# and populate the new row with this week's price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 76; this shifts rows 76..188 down to 77..189
# and carries the row formatting (e.g. the date style on column D) automatically.
$ws.Rows(76).Insert()

# Populate the newly inserted row 76 with the new record's data.
$ws.Cells.Item(76, 1).Value2 = 7
$ws.Cells.Item(76, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(76, 3).Value2 = "Ñuble"
$ws.Cells.Item(76, 4).Value2 = 45175
$ws.Cells.Item(76, 5).Value2 = 16
$ws.Cells.Item(76, 6).Value2 = "Fruta"
$ws.Cells.Item(76, 7).Value2 = 100108
$ws.Cells.Item(76, 8).Value2 = "Tropicales y subtropicales"
$ws.Cells.Item(76, 9).Value2 = 100108002
$ws.Cells.Item(76, 10).Value2 = "Mango"
$ws.Cells.Item(76, 11).Value2 = "Sin especificar"
$ws.Cells.Item(76, 12).Value2 = "Primera"
$ws.Cells.Item(76, 13).Value2 = 60
$ws.Cells.Item(76, 14).Value2 = 11000
$ws.Cells.Item(76, 15).Value2 = 11000
$ws.Cells.Item(76, 16).Value2 = 11000
$ws.Cells.Item(76, 17).Value2 = "$/bandeja 4 kilos"
$ws.Cells.Item(76, 18).Value2 = "Brasil"
$ws.Cells.Item(76, 19).Value2 = 2750
$ws.Cells.Item(76, 20).Value2 = 4
